$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Measures")
$ws2 = $wb.Worksheets.Item("Source Information")

$ws1.Range("B2").Value = '
TOTALMTD(
    SUM(FinancialData[Sales]),
    DATESMTD(''FinancialData''[Date])
)'
$ws1.Range("E2").Value = 'This calculation sums up the total sales for the current month to date (MTD) for the given financial data by understanding the differences in dates. It takes the sum of the sales figures and the corresponding dates from the financial data to determine the current month to date total.'
$ws1.Range("B3").Value = '
TOTALQTD(
    SUM(FinancialData[Sales]),
    DATESQTD(''FinancialData''[Date])
)'
$ws1.Range("E3").Value = 'This calculation totals the sales from the FinancialData table for the current quarter up to the date specified in the Date column. It is useful for getting an understanding of how the business is performing in a given quarter up to a particular point in time.'
$ws1.Range("B4").Value = '
TOTALYTD(
    SUM(FinancialData[ Sales ]),
    DATESYTD(''FinancialData''[Date])
)'
$ws1.Range("E4").Value = 'This calculation is for calculating the total year-to-date sales from a given date. It adds up all of the sales from the current year up until the specified date from the ''FinancialData'' table and stores it in a total year-to-date sales amount.'
$ws1.Range("B5").Value = '
CALCULATE(
    SUM(FinancialData[Sales]),
    PARALLELPERIOD(
        DATESMTD(FinancialData[Date].[Date]),
        -1, 
        MONTH
    )
)'
$ws1.Range("E5").Value = 'This calculation takes the sum of sales during the current month-to-date period, and then subtracts the sum of sales from the preceding period. This allows us to compare the current performance to the previous period, providing insight into the sales trend.'
$ws1.Range("B6").Value = '
CALCULATE( 
    SUM(FinancialData[Sales]), 
    DATESYTD( 
        PARALLELPERIOD(
            FinancialData[Date].[Date], 
            -1, 
            QUARTER
        ) 
    )
)'
$ws1.Range("E6").Value = 'This calculation calculates the sum of sales from the start of the same quarter of the previous year''s date, until the current date. This calculation is useful to compare year-over-year sales trends, or to calculate quarter-over-quarter sales growth.'
$ws1.Range("B7").Value = '
(Previous Month Sales MTD - SalesMTD) / Previous Month Sales MTD'
$ws1.Range("E7").Value = 'This calculation is used to find the percentage change in monthly sales from the previous month. It takes the current month''s sales MTD subtracts the sales MTD from the previous month and then divides it by the previous month''s sales MTD. This gives a percentage change of current month''s sales MTD relative'
$ws1.Range("B8").Value = '
CALCULATE (
    SUM(FinancialData[Sales]),
    DATESYTD(SAMEPERIODLASTYEAR(FinancialData[Date].[Date]))
)'
$ws1.Range("E8").Value = 'This calculation is finding the total sales value for the current year to date (YTD), compared with the same period last year (YTD). The calculation is comparing the sales data this year with the sales data from the same period last year to get an understanding of year-on-year performance.'
$ws1.Range("B9").Value = '
CALCULATE(
    SUM(FinancialData[Sales]),
    DATESBETWEEN(
        FinancialData[Date],
        DATE(2014, 1, 1),
        DATE(2014, 6, 31)
    )
)'
$ws1.Range("E9").Value = 'This calculation is summing up the total sales from the FinancialData table between the dates of January 1, 2014 and June 31, 2014. This is a useful calculation to calculate how much money was made in a specific time period.'

$ws2.Range("I2").Value = '1. "Changed Type" is a command that changes the data types of specified columns in the FinancialData_Table to text, number, date, and integer, respectively.
2. This statement renames the existing column called "COGS" to "Cost" in the table "Changed Type".
3. This changes the type of some columns in the table to currencies, dates, and integers.
4. The "Added Custom" step adds a new column to a previously changed table which calculates the total cost by adding the cost and discounts columns.
5. This sentence changes the data type of the column "Total Cost" to a currency type within the dataset "#Added Custom".

'

Write-Host "done"
